# Update water sampling logs:
# Add a new "RMI 2025 Spring Recovery" row to the SampleLog sheet, and
# resize a few columns to fit the new (wider) content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SampleLog")

# --- New row 3 values -------------------------------------------------
$ws.Range("A3").Value = "RMI"
$ws.Range("B3").Value = "RMI 2025 Spring Recovery"
$ws.Range("C3").Value = "ru39-20250423T1535"
$ws.Range("D3").Value = "recovery"

# Date UTC (5/20/2025) - copy formatting from the cell above (E2) so the
# new cell picks up the existing date number format, then set the value.
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = 45797

# Deployment Notes - copy formatting from the cell above (F2) so the new
# cell picks up the existing wrap-text format, then set the long note.
$notes = "Collected extra deep samples, 4 bottles (2 pH, 2 TA) from 2 niskins per cast.`nCast2 didn't close bottles, so redid cast`nRosette didn't close 2m bottle on 3nd cast, so collected 2m sample immediately after rosette back on board for the 3nd profile samples.`nThere is some issue with the ctd talking to the rosette, maybe a connection or cable issue? The bottle files had jumbled/weird characters, so couldn't process those files. Temp/Sal are from the CTD binned cnv files from the ctd on the rosette during the water collection. People: Nicole Waite, Brian Buckingham, Jess Leonard and Delphine Mossman"
$ws.Range("F2").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = $notes

# New row is a multi-line note, so it needs to be taller than the default.
$ws.Rows.Item(3).RowHeight = 128

# --- Column widths ------------------------------------------------------
# The new row's longer entries mean a few columns need to be a bit wider
# to fit (mirrors Excel's own "best fit" column auto-sizing).
$ws.Columns.Item(2).ColumnWidth = 25.998697916666668
$ws.Columns.Item(3).ColumnWidth = 17.998697916666668
$ws.Columns.Item(4).ColumnWidth = 22.498697916666668
$ws.Columns.Item(5).ColumnWidth = 8.830729166666666
